# Auto-generated: update TPM-derived values in Efnb1-Ephb2 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.546140333333334
$ws.Range("H2").Value = 28.638421
$ws.Range("I2").Value = 0.587227294878132
$ws.Range("J2").Value = 0.587227294878132
$ws.Range("O2").Value = 0.9347132976570145
$ws.Range("P2").Value = 0.9347132976570145
$ws.Range("Q2").Value = 84.79579030515677
$ws.Range("R2").Value = 763.1621127464109
$ws.Range("S2").Value = 0.5488891612697467
$ws.Range("T2").Value = 0.5488891612697467

# Row 3
$ws.Range("G3").Value = 9.546140333333334
$ws.Range("H3").Value = 28.638421
$ws.Range("I3").Value = 0.587227294878132
$ws.Range("J3").Value = 0.587227294878132
$ws.Range("M3").Value = 0.616144
$ws.Range("N3").Value = 1.848432
$ws.Range("O3").Value = 0.06483569448352988
$ws.Range("P3").Value = 0.0648356944835299
$ws.Range("Q3").Value = 5.881797089541334
$ws.Range("R3").Value = 52.936173805872
$ws.Range("S3").Value = 0.03807328948310828
$ws.Range("T3").Value = 0.03807328948310829

# Row 4
$ws.Range("G4").Value = 9.546140333333334
$ws.Range("H4").Value = 28.638421
$ws.Range("I4").Value = 0.587227294878132
$ws.Range("J4").Value = 0.587227294878132
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.004286
$ws.Range("N4").Value = 0.012858
$ws.Range("O4").Value = 0.0004510078594555965
$ws.Range("P4").Value = 0.0004510078594555965
$ws.Range("Q4").Value = 0.04091475746866666
$ws.Range("R4").Value = 0.368232817218
$ws.Range("S4").Value = 0.0002648441252768867
$ws.Range("T4").Value = 0.0002648441252768867

# Row 5
$ws.Range("I5").Value = 0.2496684258894083
$ws.Range("J5").Value = 0.2496684258894083
$ws.Range("O5").Value = 0.9347132976570145
$ws.Range("P5").Value = 0.9347132976570145
$ws.Range("S5").Value = 0.2333683976839247
$ws.Range("T5").Value = 0.2333683976839247

# Row 6
$ws.Range("I6").Value = 0.2496684258894083
$ws.Range("J6").Value = 0.2496684258894083
$ws.Range("M6").Value = 0.616144
$ws.Range("N6").Value = 1.848432
$ws.Range("O6").Value = 0.06483569448352988
$ws.Range("P6").Value = 0.0648356944835299
$ws.Range("Q6").Value = 2.500733589114667
$ws.Range("R6").Value = 22.506602302032
$ws.Range("S6").Value = 0.0161874257831495
$ws.Range("T6").Value = 0.0161874257831495

# Row 7
$ws.Range("I7").Value = 0.2496684258894083
$ws.Range("J7").Value = 0.2496684258894083
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.004286
$ws.Range("N7").Value = 0.012858
$ws.Range("O7").Value = 0.0004510078594555965
$ws.Range("P7").Value = 0.0004510078594555965
$ws.Range("Q7").Value = 0.01739551819533334
$ws.Range("R7").Value = 0.156559663758
$ws.Range("S7").Value = 0.0001126024223340303
$ws.Range("T7").Value = 0.0001126024223340303

# Row 8
$ws.Range("G8").Value = 2.210442
$ws.Range("H8").Value = 6.631326
$ws.Range("I8").Value = 0.1359745227725727
$ws.Range("J8").Value = 0.1359745227725727
$ws.Range("O8").Value = 0.9347132976570145
$ws.Range("P8").Value = 0.9347132976570145
$ws.Range("Q8").Value = 19.634760203474
$ws.Range("R8").Value = 176.712841831266
$ws.Range("S8").Value = 0.1270971945780902
$ws.Range("T8").Value = 0.1270971945780902

# Row 9
$ws.Range("G9").Value = 2.210442
$ws.Range("H9").Value = 6.631326
$ws.Range("I9").Value = 0.1359745227725727
$ws.Range("J9").Value = 0.1359745227725727
$ws.Range("M9").Value = 0.616144
$ws.Range("N9").Value = 1.848432
$ws.Range("O9").Value = 0.06483569448352988
$ws.Range("P9").Value = 0.0648356944835299
$ws.Range("Q9").Value = 1.361950575648
$ws.Range("R9").Value = 12.257555180832
$ws.Range("S9").Value = 0.008816002616026298
$ws.Range("T9").Value = 0.0088160026160263

# Row 10
$ws.Range("G10").Value = 2.210442
$ws.Range("H10").Value = 6.631326
$ws.Range("I10").Value = 0.1359745227725727
$ws.Range("J10").Value = 0.1359745227725727
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.004286
$ws.Range("N10").Value = 0.012858
$ws.Range("O10").Value = 0.0004510078594555965
$ws.Range("P10").Value = 0.0004510078594555965
$ws.Range("Q10").Value = 0.009473954412
$ws.Range("R10").Value = 0.08526558970800001
$ws.Range("S10").Value = 0.00006132557845615426
$ws.Range("T10").Value = 0.00006132557845615426

# Row 11
$ws.Range("G11").Value = 0.4410293333333333
$ws.Range("H11").Value = 1.323088
$ws.Range("I11").Value = 0.02712975645988715
$ws.Range("J11").Value = 0.02712975645988715
$ws.Range("O11").Value = 0.9347132976570145
$ws.Range("P11").Value = 0.9347132976570145
$ws.Range("Q11").Value = 3.917544637089778
$ws.Range("R11").Value = 35.257901733808
$ws.Range("S11").Value = 0.02535854412525281
$ws.Range("T11").Value = 0.02535854412525281

# Row 12
$ws.Range("G12").Value = 0.4410293333333333
$ws.Range("H12").Value = 1.323088
$ws.Range("I12").Value = 0.02712975645988715
$ws.Range("J12").Value = 0.02712975645988715
$ws.Range("M12").Value = 0.616144
$ws.Range("N12").Value = 1.848432
$ws.Range("O12").Value = 0.06483569448352988
$ws.Range("P12").Value = 0.0648356944835299
$ws.Range("Q12").Value = 0.2717375775573334
$ws.Range("R12").Value = 2.445638198016
$ws.Range("S12").Value = 0.001758976601245814
$ws.Range("T12").Value = 0.001758976601245815

# Row 13
$ws.Range("G13").Value = 0.4410293333333333
$ws.Range("H13").Value = 1.323088
$ws.Range("I13").Value = 0.02712975645988715
$ws.Range("J13").Value = 0.02712975645988715
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.004286
$ws.Range("N13").Value = 0.012858
$ws.Range("O13").Value = 0.0004510078594555965
$ws.Range("P13").Value = 0.0004510078594555965
$ws.Range("Q13").Value = 0.001890251722666667
$ws.Range("R13").Value = 0.017012265504
$ws.Range("S13").Value = 0.00001223573338852534
$ws.Range("T13").Value = 0.00001223573338852535

Write-Output "Updated cells"